$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2 through 81 from serial date 45185 to 45204
for ($row = 2; $row -le 81; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45185) {
        $cell.Value2 = 45204
    }
}
